# Add new fields (columns) to the "Tabelle5" table on the
# "Booklet_FK Lagerlogistik" sheet: AssessmentType, Description, Disclaimer,
# Duration, EscoOccupationId, EscoSkills, Publisher, Title.
# New column "AssessmentType" is populated with 0 for every existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Booklet_FK Lagerlogistik")
$ws.Activate()

$lo = $ws.ListObjects.Item(1)

$newColumnNames = @("AssessmentType", "Description", "Disclaimer", "Duration", "EscoOccupationId", "EscoSkills", "Publisher", "Title")

# Remember the formatting of the last existing header cell so the new
# headers match the look of the rest of the header row.
$lastHeaderCell = $lo.HeaderRowRange.Cells.Item(1, $lo.HeaderRowRange.Columns.Count)
$lastHeaderCell.Copy()

foreach ($name in $newColumnNames) {
    $col = $lo.ListColumns.Add()
    $headerCell = $col.Range.Cells.Item(1, 1)
    $headerCell.Value = $name
    $headerCell.PasteSpecial(-4122)
}

# Populate the new "AssessmentType" column with 0 for every data row.
$assessmentTypeCol = $lo.ListColumns.Item("AssessmentType")
$dataRange = $assessmentTypeCol.DataBodyRange
for ($i = 1; $i -le $dataRange.Rows.Count; $i++) {
    $dataRange.Cells.Item($i, 1).Value = 0
}

# Match the author's final selection after adding the columns.
$ws.Range("BJ2").Select()
